$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing last data row (214), pushing it
# down to row 216. This also carries the existing column formatting
# (e.g. the date style on column D) into the new rows.
$ws.Range("A214:R215").EntireRow.Insert()

# New row 214: weekly "Primera" quality record for 2022-02-18
$ws.Range("A214").Value = 11
$ws.Range("B214").Value = "Vega Monumental Concepción"
$ws.Range("C214").Value = "Bíobío"
$ws.Range("D214").Value = 44610
$ws.Range("E214").Value = 8
$ws.Range("F214").Value = 100112008
$ws.Range("G214").Value = "Coliflor"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 2000
$ws.Range("K214").Value = 1000
$ws.Range("L214").Value = 1000
$ws.Range("M214").Value = 1000
$ws.Range("N214").Value = "$/unidad"
$ws.Range("O214").Value = "Región Metropolitana"
$ws.Range("P214").Value = 1000
$ws.Range("Q214").Value = 1
$ws.Range("R214").Value = "Hortaliza"

# New row 215: weekly "Segunda" quality record for 2022-02-18
$ws.Range("A215").Value = 11
$ws.Range("B215").Value = "Vega Monumental Concepción"
$ws.Range("C215").Value = "Bíobío"
$ws.Range("D215").Value = 44610
$ws.Range("E215").Value = 8
$ws.Range("F215").Value = 100112008
$ws.Range("G215").Value = "Coliflor"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Segunda"
$ws.Range("J215").Value = 1000
$ws.Range("K215").Value = 800
$ws.Range("L215").Value = 800
$ws.Range("M215").Value = 800
$ws.Range("N215").Value = "$/unidad"
$ws.Range("O215").Value = "Región Metropolitana"
$ws.Range("P215").Value = 800
$ws.Range("Q215").Value = 1
$ws.Range("R215").Value = "Hortaliza"
